$d = $word.ActiveDocument
$rng = $d.Paragraphs.Last.Range
$xml = '<w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Les </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>tranformations</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  &#8230;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Titre3"/></w:pPr><w:r><w:t xml:space="preserve">Les styles </w:t></w:r></w:p><w:p><w:r><w:t>D&#233;finissent l&#8217;apparence des contr&#244;les.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Utilisation de Setters qui d&#233;finissent la propri&#233;t&#233; &#224; modifier gr&#226;ce &#224; la propri&#233;t&#233; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Property</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> et la valeur (propri&#233;t&#233; Value).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Titre3"/></w:pPr><w:r><w:t xml:space="preserve">Les </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>templates</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:t>D&#233;crivent la structure visuelle d&#8217;un contr&#244;le. Propri&#233;t&#233; Template.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Titre3"/></w:pPr><w:r><w:t>Les triggers</w:t></w:r></w:p><w:p><w:r><w:t>Utiliser pour r&#233;aliser des applications proposant des interactions riches et dynamiques.</w:t></w:r></w:p><w:p><w:r><w:t>WPF v&#233;rifie trois choses dans les conditions d&#8217;un trigger&#160;:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Une </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Property</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dependency</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (Trigger)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr><w:r><w:t>Propri&#233;t&#233; .NET (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DataTrigger</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="21"/></w:numPr></w:pPr><w:r><w:t>Ev&#233;nement (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EventTrigger</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>Dans les deux premi&#232;res conditions, le trigger est d&#233;clench&#233; lorsque la propri&#233;t&#233; sp&#233;cifi&#233;e est modifi&#233;e.</w:t></w:r></w:p><w:p><w:r><w:t>Dans le dernier cas, d&#233;clenchement lorsque l&#8217;&#233;v&#233;nement survient.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Pour modifier la valeur d&#8217;une propri&#233;t&#233; dans un trigger, on utilise Setter.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Titre3"/></w:pPr><w:r><w:t>Les storyboards</w:t></w:r></w:p><w:p><w:r><w:t>El&#233;ments XAML qui permettent de d&#233;finir un ensemble d&#8217;actions. C&#8217;est un ensemble d&#8217;animations/transformations. Ils permettent un param&#233;trage complet des animations.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Titre3"/></w:pPr><w:r><w:t>Les animations / transformations</w:t></w:r></w:p><w:p><w:r><w:t>Int&#233;ressant mais pas utiles pour le projet normalement</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Titre3"/></w:pPr><w:r><w:t>3D</w:t></w:r></w:p><w:p><w:r><w:t>Pas utile pour notre projet</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$rng.InsertXML($xml)
